# Rename "Sheet1" -> "Fig.5" and move the active selection to E17,
# matching the 2nd-revision package update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Fig.5"
$ws.Activate()
$ws.Range("E17").Select()
